# Refine skills section to showcase modern web deployment expertise
$d = $word.ActiveDocument

# 1. Update the section heading
$d.Content.Find.Execute(
    "Web Development", $true, $false, $false, $false, $false,
    $true, 1, $false, "No-Code / Low-Code Web Deployment", 2
)

# 2. Replace the four bullet lines under the section
$d.Content.Find.Execute(
    "HTML5, CSS3, JavaScript (Responsive Design, Interactive UIs)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Replit: Web Page Design & Deployment (HTML Structuring, Visual Flow, Script Integration)", 2
)

$d.Content.Find.Execute(
    "Node.js, Express.js, TypeScript (Backend Development, API Integration)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Layout & Usability Architecture (Fonts, Spacing, Mobile-Friendliness, User Interaction)", 2
)

$d.Content.Find.Execute(
    "PostgreSQL, Drizzle ORM (Database Management, Full-Stack Development)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Integration of Tools: Formspree, Netlify, Cloudflare for Functional Delivery", 2
)

$d.Content.Find.Execute(
    "Replit Development, Static Hosting, GitHub Pages Deployment", $true, $false, $false, $false, $false,
    $true, 1, $false, "Full-stack launch using prebuilt components and code remixing", 2
)

# 3. Insert a brand-new bullet paragraph right after the "Full-stack launch..." bullet,
#    inheriting the same paragraph formatting (spacing after = 60)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Full-stack launch using prebuilt components and code remixing*") {
        $newPara = $p.Range.InsertParagraphAfter()
        $p.Next().Range.Text = "• UX-Focused Page Building: Professional design without hand-coding"
        break
    }
}
